# Future Cart title-slide relayout: reposition/resize the title, subtitle
# and "Presented By" text boxes, and bump several font sizes / typefaces.
#
# NOTE on units: Shape.Left/Top/Width/Height are expressed in points while
# the underlying OOXML stores EMU (1 pt = 12700 EMU). The host round-trips
# the point value through a 32-bit float before truncating to EMU, so a
# naive `emu/12700.0` can land one EMU short of the intended value. EmuToPt
# binary-searches the smallest double that, once rounded to single
# precision and multiplied back out, truncates to exactly the requested
# EMU count - this keeps the shape geometry pixel-exact vs. the source
# OOXML.
function EmuToPt {
    param([double]$Emu)
    $loF = $Emu / 12700.0
    $hiF = ($Emu + 1) / 12700.0
    $loPt = $loF - 0.01
    $hiPt = $hiF + 0.01
    for ($i = 0; $i -lt 100; $i++) {
        $mid = ($loPt + $hiPt) / 2.0
        $midAsSingle = [double]([single]$mid)
        if ($midAsSingle -ge $loF) {
            $hiPt = $mid
        } else {
            $loPt = $mid
        }
    }
    return $hiPt
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---- Title 1 : just slides up/left a bit ----------------------------
$title = $s.Shapes.Item(1)
$title.Left = EmuToPt 549110
$title.Top  = EmuToPt 513662

# ---- Subtitle 2 : reposition/resize + bigger subtitle text ----------
$subtitle = $s.Shapes.Item(2)
$subtitle.Left   = EmuToPt 2601797
$subtitle.Top    = EmuToPt 1868046
$subtitle.Width  = EmuToPt 9907572
$subtitle.Height = EmuToPt 685800

$subtitleRange = $subtitle.TextFrame.TextRange.Characters(1, 60)
$subtitleRange.Font.Size = 28
$subtitleRange.Font.NameAscii = "Calibri"
$subtitleRange.Font.NameFarEast = "Calibri"
$subtitleRange.Font.NameComplexScript = "Calibri"

# ---- TextBox 5 ("Presented By : ...") : reposition/resize + restyle -
$credit = $s.Shapes.Item(3)
$credit.Left   = EmuToPt 3158936
$credit.Top    = EmuToPt 4172180
$credit.Width  = EmuToPt 5111271
$credit.Height = EmuToPt 1238801

$creditRange = $credit.TextFrame.TextRange

# "Presented By : " run -> 32pt, Aptos
$presentedByRun = $creditRange.Characters(1, 15)
$presentedByRun.Font.Size = 32
$presentedByRun.Font.NameAscii = "Aptos"

# leading spaces run -> italic, Calibri
$spacesRun = $creditRange.Characters(18, 13)
$spacesRun.Font.Italic = $true
$spacesRun.Font.NameAscii = "Calibri"
$spacesRun.Font.NameFarEast = "Calibri"
$spacesRun.Font.NameComplexScript = "Calibri"

# name run -> 32pt, italic, Calibri
$nameRun = $creditRange.Characters(31, 23)
$nameRun.Font.Size = 32
$nameRun.Font.Italic = $true
$nameRun.Font.NameAscii = "Calibri"
$nameRun.Font.NameFarEast = "Calibri"
$nameRun.Font.NameComplexScript = "Calibri"
